$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. The former 5th shared string "/conclusions/approve1" becomes "Assert400"
#    and is re-purposed as a new section header (like "Assert200"/"Assert401").
#    Row 3 (which used to display "/conclusions/approve1") now shows
#    "/conclusions/approve" instead (same text already used by row 11), and
#    its row shrinks to the shorter wrapped height.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = "/conclusions/approve"
$ws.Range("A3").RowHeight = 30

# ---------------------------------------------------------------------------
# 2. Add a brand new "Assert400" test-case block in rows 14-16, mirroring the
#    existing "Assert401" block (rows 9-11): a merged header row, an
#    "EndPoint" label row, and a path value row.
# ---------------------------------------------------------------------------
$ws.Range("A9:E11").Copy()
$ws.Range("A14:E16").PasteSpecial(-4122)

$ws.Range("A14").Value2 = "Assert400"
$ws.Range("A15").Value2 = "EndPoint"
$ws.Range("A16").Value2 = "/conclusions/approve"
$ws.Range("A16").RowHeight = 30

$ws.Range("A14:E14").Merge()

# ---------------------------------------------------------------------------
# 3. Update the active selection to reflect where the editor left off.
# ---------------------------------------------------------------------------
$ws.Range("G24").Select()
